$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Replace numeric values in D2:D5 with text values "A123".."A126"
$ws.Range("D2").Value = "A123"
$ws.Range("D3").Value = "A124"
$ws.Range("D4").Value = "A125"
$ws.Range("D5").Value = "A126"

# Update the selection to match the new active cell / selection range
$ws.Range("D2:D5").Select()
